$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append the new data row (year 2024) right after the last existing row (25),
# copying the formatting (borders/alignment) used by the row above it.
$ws.Range("A25:D25").Copy()
$ws.Range("A26:D26").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("A26").Value = 2024
$ws.Range("B26").Value = 7
$ws.Range("C26").Value = 39
$ws.Range("D26").Value = 7.14

# Update the view: scroll so row 13 is near the top, and move the selection
# to reflect where the user ended up after entering the new row.
$excel.Windows.Item(1).ScrollRow = 13
$ws.Range("B29").Select()
